# Auto-generated edit script applying the diff to before.xlsx
# Updates "想去人数" (F col) counts across sheets, one G-column
# availability cell, and refreshes two rows in "全部类型" so they
# match the corresponding (now-updated) rows in "演出".

$wb = $excel.ActiveWorkbook

# ---- Sheet "展览" (sheet1) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value2 = 88
$ws1.Range("F4").Value2 = 184
$ws1.Range("F5").Value2 = 996
$ws1.Range("F7").Value2 = 2485
$ws1.Range("F9").Value2 = 1215
$ws1.Range("F10").Value2 = 892
$ws1.Range("F12").Value2 = 897
$ws1.Range("F13").Value2 = 1102
$ws1.Range("F15").Value2 = 290
$ws1.Range("F17").Value2 = 719
$ws1.Range("F18").Value2 = 754
$ws1.Range("F19").Value2 = 183
$ws1.Range("F20").Value2 = 473
$ws1.Range("F21").Value2 = 1102
$ws1.Range("F22").Value2 = 78
$ws1.Range("F23").Value2 = 572
$ws1.Range("F24").Value2 = 580
$ws1.Range("F25").Value2 = 214
$ws1.Range("F26").Value2 = 294
$ws1.Range("F27").Value2 = 291
$ws1.Range("F28").Value2 = 676
$ws1.Range("F29").Value2 = 85
$ws1.Range("F30").Value2 = 3786
$ws1.Range("F31").Value2 = 474
$ws1.Range("F32").Value2 = 42
$ws1.Range("F36").Value2 = 134
$ws1.Range("F37").Value2 = 1583
$ws1.Range("F38").Value2 = 436
$ws1.Range("F39").Value2 = 142
$ws1.Range("F41").Value2 = 134
$ws1.Range("F42").Value2 = 65
$ws1.Range("F44").Value2 = 122
$ws1.Range("F46").Value2 = 88

# ---- Sheet "演出" (sheet2) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F8").Value2 = 7
$ws2.Range("F11").Value2 = 177
$ws2.Range("F13").Value2 = 5
$ws2.Range("F16").Value2 = 176

# ---- Sheet "本地生活" (sheet3) ----
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value2 = 2254
$ws3.Range("F3").Value2 = 717
$ws3.Range("F4").Value2 = 664
$ws3.Range("G4").Value2 = "不可售"  # was numeric 20, now text status

# ---- Sheet "全部类型" (sheet4) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value2 = 2254
$ws4.Range("F3").Value2 = 717
$ws4.Range("F7").Value2 = 996
$ws4.Range("F8").Value2 = 2485
$ws4.Range("F10").Value2 = 1215
$ws4.Range("F11").Value2 = 892
$ws4.Range("F13").Value2 = 897
$ws4.Range("F14").Value2 = 1102
$ws4.Range("F15").Value2 = 290
$ws4.Range("F18").Value2 = 719
$ws4.Range("F21").Value2 = 754
$ws4.Range("F22").Value2 = 183
$ws4.Range("F23").Value2 = 473
$ws4.Range("F24").Value2 = 1102
$ws4.Range("F26").Value2 = 78
$ws4.Range("F27").Value2 = 572
$ws4.Range("F28").Value2 = 580
$ws4.Range("F29").Value2 = 214
$ws4.Range("F30").Value2 = 294
$ws4.Range("F31").Value2 = 291
$ws4.Range("F33").Value2 = 3789
$ws4.Range("F34").Value2 = 177
$ws4.Range("F35").Value2 = 474
$ws4.Range("F36").Value2 = 42
$ws4.Range("F38").Value2 = 134
$ws4.Range("F39").Value2 = 1583
$ws4.Range("F40").Value2 = 436
$ws4.Range("F41").Value2 = 142
$ws4.Range("F44").Value2 = 134
$ws4.Range("F46").Value2 = 122
$ws4.Range("F47").Value2 = 88

# Row 5 in "全部类型" is replaced wholesale with the event that used
# to be row 5 in "演出" (井草圣二 指弹吉他音乐会).
# B5 holds a plain "YYYY-MM-DD" label that must stay TEXT, not get
# auto-converted to a date serial by the COM date parser, so force
# the cell to Text format before assigning it.
$ws4.Range("B5").NumberFormat = "@"
$ws4.Range("B5").Value2 = "2024-08-02"
$ws4.Range("C5").Value2 = "广州·井草圣二 2024《夏日独白》指弹吉他音乐会"
$ws4.Range("D5").Value2 = "恩宁路265号3层 MaoLivehouse(永庆坊店)"
$ws4.Range("E5").Value2 = "2024.08.02 19:30-08.02 21:00"
$ws4.Range("F5").Value2 = 19
$ws4.Range("G5").Value2 = 260
$ws4.Range("H5").Value2 = "https://show.bilibili.com/platform/detail.html?id=86940"
$ws4.Range("I5").Value2 = "//i0.hdslb.com/bfs/openplatform/202406/iNGVydXM1717644835981.jpeg"

# Row 6 in "全部类型" is replaced wholesale with the event that used
# to be row 6 in "演出" (忱宴·渐渐被你吸引 ACG演唱会).
$ws4.Range("B6").NumberFormat = "@"
$ws4.Range("B6").Value2 = "2024-08-03"
$ws4.Range("C6").Value2 = "广州·【暑期5折】《忱宴·渐渐被你吸引》热血动漫二次元ACG演唱会"
$ws4.Range("D6").Value2 = "东风中路299号 广州中山纪念堂"
$ws4.Range("E6").Value2 = "2024.08.03 20:00-08.03 21:40"
$ws4.Range("F6").Value2 = 80
$ws4.Range("G6").Value2 = 50
$ws4.Range("H6").Value2 = "https://show.bilibili.com/platform/detail.html?id=85917"
$ws4.Range("I6").Value2 = "//i1.hdslb.com/bfs/openplatform/202405/won43hte1715675570347.jpeg"
